# Atualizacao da programacao da disciplina
# Fill in rows 8-12 of the schedule table with the new "CyberBattlesim"
# group (4 students) and the start of the "Markov Decision Process" group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: new theme/title row (CyberBattlesim), first student = Davi ---
$ws.Range("A8").Value = "Ambiente CyberBattlesim com Reinforcement Learning"
$ws.Range("B8").Value = "Davi"
$ws.Range("C8").Value = "Entregue"

$ws.Range("A8").WrapText = $true
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Style = "Good"
$ws.Range("C8").WrapText = $true

$ws.Range("A8:F8").RowHeight = 28.5

# --- Row 9: Henrique ---
$ws.Range("B9").Value = "Henrique"
$ws.Range("C9").Value = "Entregue"

$ws.Range("B9").WrapText = $true
$ws.Range("C9").Style = "Good"
$ws.Range("C9").WrapText = $true

# --- Row 10: Luiza ---
$ws.Range("B10").Value = "Luiza"
$ws.Range("C10").Value = "Entregue"

$ws.Range("B10").WrapText = $true
$ws.Range("C10").Style = "Good"
$ws.Range("C10").WrapText = $true

# --- Row 11: Nicolas ---
$ws.Range("B11").Value = "Nicolas"
$ws.Range("C11").Value = "Entregue"

$ws.Range("B11").WrapText = $true
$ws.Range("C11").Style = "Good"
$ws.Range("C11").WrapText = $true

# --- Row 12: new theme/title row (Markov Decision Process), first student = Ana Carolina ---
$ws.Range("A12").Value = "Markov Decision Process no ambiente SimpleGrid do Gymnasium"
$ws.Range("B12").Value = "Ana Carolina"
$ws.Range("C12").Value = "entregue"

$ws.Range("A12").WrapText = $true
$ws.Range("B12").WrapText = $true
$ws.Range("C12").Style = "Good"
$ws.Range("C12").WrapText = $true

$ws.Range("A12:F12").RowHeight = 42.75
